{"js": "// (find text, replacement text) pairs -- one per run that changed in this\n// revision: the header date plus each \"NNN\u00d7N=\" multiplication prompt in\n// the table.\nconst pairs = [\n  [\"2026-02-21 Saturday\", \"2026-02-22 Sunday\"],\n  [\"372\u00d77=\", \"811\u00d79=\"],\n  [\"326\u00d72=\", \"576\u00d78=\"],\n  [\"130\u00d79=\", \"411\u00d79=\"],\n  [\"450\u00d78=\", \"159\u00d72=\"],\n  [\"983\u00d78=\", \"515\u00d77=\"],\n  [\"799\u00d79=\", \"379\u00d77=\"],\n  [\"586\u00d73=\", \"411\u00d78=\"],\n  [\"109\u00d75=\", \"275\u00d78=\"],\n  [\"498\u00d77=\", \"778\u00d76=\"],\n  [\"247\u00d79=\", \"157\u00d75=\"],\n  [\"459\u00d78=\", \"573\u00d79=\"],\n  [\"953\u00d75=\", \"806\u00d76=\"],\n  [\"547\u00d73=\", \"987\u00d75=\"],\n  [\"875\u00d75=\", \"591\u00d77=\"],\n  [\"607\u00d78=\", \"142\u00d74=\"],\n  [\"921\u00d78=\", \"434\u00d74=\"],\n  [\"630\u00d79=\", \"850\u00d77=\"],\n  [\"956\u00d77=\", \"149\u00d77=\"],\n  [\"267\u00d72=\", \"165\u00d78=\"],\n  [\"718\u00d76=\", \"212\u00d79=\"],\n  [\"242\u00d78=\", \"331\u00d77=\"],\n  [\"545\u00d77=\", \"864\u00d75=\"],\n  [\"930\u00d77=\", \"616\u00d79=\"],\n  [\"124\u00d77=\", \"345\u00d75=\"],\n  [\"332\u00d72=\", \"943\u00d79=\"],\n];\n\nconst body = context.document.body;\nconst allResults = [];\n\nfor (const [findText] of pairs) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  allResults.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < pairs.length; i++) {\n  const [, replaceText] = pairs[i];\n  const results = allResults[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(replaceText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# (find text, replacement text) pairs -- one per run that changed in this revision:\n# the header date plus each \"NNN\u00d7N=\" multiplication prompt in the table.\n$pairs = @(\n    @(\"2026-02-21 Saturday\", \"2026-02-22 Sunday\"),\n    @(\"372\u00d77=\", \"811\u00d79=\"),\n    @(\"326\u00d72=\", \"576\u00d78=\"),\n    @(\"130\u00d79=\", \"411\u00d79=\"),\n    @(\"450\u00d78=\", \"159\u00d72=\"),\n    @(\"983\u00d78=\", \"515\u00d77=\"),\n    @(\"799\u00d79=\", \"379\u00d77=\"),\n    @(\"586\u00d73=\", \"411\u00d78=\"),\n    @(\"109\u00d75=\", \"275\u00d78=\"),\n    @(\"498\u00d77=\", \"778\u00d76=\"),\n    @(\"247\u00d79=\", \"157\u00d75=\"),\n    @(\"459\u00d78=\", \"573\u00d79=\"),\n    @(\"953\u00d75=\", \"806\u00d76=\"),\n    @(\"547\u00d73=\", \"987\u00d75=\"),\n    @(\"875\u00d75=\", \"591\u00d77=\"),\n    @(\"607\u00d78=\", \"142\u00d74=\"),\n    @(\"921\u00d78=\", \"434\u00d74=\"),\n    @(\"630\u00d79=\", \"850\u00d77=\"),\n    @(\"956\u00d77=\", \"149\u00d77=\"),\n    @(\"267\u00d72=\", \"165\u00d78=\"),\n    @(\"718\u00d76=\", \"212\u00d79=\"),\n    @(\"242\u00d78=\", \"331\u00d77=\"),\n    @(\"545\u00d77=\", \"864\u00d75=\"),\n    @(\"930\u00d77=\", \"616\u00d79=\"),\n    @(\"124\u00d77=\", \"345\u00d75=\"),\n    @(\"332\u00d72=\", \"943\u00d79=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
